$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 83 and 84: swap all data (columns B..AC) between the two rows.
#    Column A (the running index, 81 / 82) stays put on its own row.
# ---------------------------------------------------------------------------
$row83 = @{
    B = 6227884
    C = "Canada Premier League"
    D = "Canada Premier League"
    E = 45206.75
    F = "Cavalry FC"
    G = "Pacific FC CA"
    H = 3
    I = 0
    J = "H"
    K = 2.25
    L = 3.1
    M = 2.875
    N = 2.05
    O = 3.2
    P = 3.2
    Q = -0.25
    R = 1.825
    S = 1.975
    T = 2.5
    U = 1.825
    V = 1.975
    W = 1.05
    X = -1
    Y = -1
    Z = 0.825
    AA = -1
    AB = 0.825
    AC = -1
}

$row84 = @{
    B = 7301364
    C = "Canada Premier League"
    D = "Canada Premier League"
    E = 45206.75
    F = "Forge FC"
    G = "Atletico Ottawa"
    H = 0
    I = 1
    J = "A"
    K = 1.8
    L = 3.6
    M = 3.5
    N = 1.533
    O = 3.8
    P = 5
    Q = -1
    R = 1.975
    S = 1.825
    T = 2.5
    U = 1.9
    V = 1.9
    W = -1
    X = -1
    Y = 4
    Z = -1
    AA = 0.825
    AB = -1
    AC = 0.8999999999999999
}

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($col in $cols) {
    $ws.Range("$col" + "83").Value2 = $row83[$col]
    $ws.Range("$col" + "84").Value2 = $row84[$col]
}

# ---------------------------------------------------------------------------
# 2) A new match result comes in at row 93 (91st game), pushing the former
#    row 93 (still unplayed at the time) down to row 94.
# ---------------------------------------------------------------------------

# Move the formatting of the two styled cells in row 93 (A: bold/border id
# style, E: custom date style) down to row 94 before the values travel, so
# row 94 ends up dressed exactly like row 93 used to be - without inserting
# a whole row (which would pull in unrelated formatting / blow out the
# dimension).
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E93").Copy() | Out-Null
$ws.Range("E94").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Move row 93's current values down into row 94, with refreshed odds.
$movedCols = @("A","B","C","D","E","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($col in $movedCols) {
    $ws.Range("$col" + "94").Value2 = $ws.Range("$col" + "93").Value2
}
$ws.Range("A94").Value2 = 92
$ws.Range("N94").Value2 = 1.65
$ws.Range("O94").Value2 = 3.8
$ws.Range("P94").Value2 = 4
$ws.Range("Q94").Value2 = -0.75
$ws.Range("R94").Value2 = 1.85
$ws.Range("S94").Value2 = 1.95
$ws.Range("U94").Value2 = 1.95
$ws.Range("V94").Value2 = 1.85

# Now overwrite row 93 with the brand-new match entry.
$newRow93 = @{
    A = 91
    B = 7803363
    C = "Canada Premier League"
    D = "Canada Premier League"
    E = 45400.95833333334
    F = "Vancouver FC"
    G = "HFX Wanderers"
    H = 2
    I = 0
    J = "H"
    K = 3.1
    L = 3.4
    M = 2
    N = 3
    O = 3.4
    P = 2.05
    Q = 0.25
    R = 1.95
    S = 1.85
    T = 2.25
    U = 1.8
    V = 2
    W = 2
    X = -1
    Y = -1
    Z = 0.95
    AA = -1
    AB = -0.5
    AC = 0.5
}

foreach ($col in $newRow93.Keys) {
    $ws.Range("$col" + "93").Value2 = $newRow93[$col]
}
